$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test run label in B2 and B7, dropping the trailing date stamp.
$ws.Range("B2").Value = "LIVEHTA Automation - Test_NonOncology_Automation_2"
$ws.Range("B7").Value = "LIVEHTA Automation - Test_NonOncology_Automation_2"

# Column B shrinks (bestFit) now that the text is shorter.
$ws.Columns.Item(2).ColumnWidth = 48

# Update the view state: selection moved to C10, scroll reset (no frozen-pane offset).
$ws.Activate()
$ws.Range("C10").Select() | Out-Null
